$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: duplicate of row 30's layout, with a new "X" marker in column F ---
$ws.Range("B30").Copy()
$ws.Range("B32").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B32").Formula = "=COUNTA(C30:I30)"

$ws.Range("C30").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$ws.Range("C32").Value = "["

$ws.Range("D30").Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("D32").Value = "&"

$ws.Range("E30").Copy()
$ws.Range("E32").PasteSpecial(-4122)
$ws.Range("E32").Value = "P"

$ws.Range("F30").Copy()
$ws.Range("F32").PasteSpecial(-4122)
$ws.Range("F32").Value = "X"

$ws.Range("G30").Copy()
$ws.Range("G32").PasteSpecial(-4122)
$ws.Range("G32").Value = "Crc1"

$ws.Range("H30").Copy()
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("H32").Value = "crc2"

$ws.Range("I30").Copy()
$ws.Range("I32").PasteSpecial(-4122)
$ws.Range("I32").Value = "]"

$ws.Range("J30").Copy()
$ws.Range("J32").PasteSpecial(-4122)
$ws.Range("J32").Value = "host - node"

# --- Row 33: same layout again, with a lowercase "x" marker, empty source row (31) ---
$ws.Range("B30").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$ws.Range("B33").Formula = "=COUNTA(C31:I31)"

$ws.Range("C30").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C33").Value = "["

$ws.Range("D30").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value = "&"

$ws.Range("E30").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value = "P"

$ws.Range("F30").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("F33").Value = "x"

$ws.Range("G30").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("G33").Value = "Crc1"

$ws.Range("H30").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("H33").Value = "crc2"

$ws.Range("I30").Copy()
$ws.Range("I33").PasteSpecial(-4122)
$ws.Range("I33").Value = "]"

$ws.Range("J30").Copy()
$ws.Range("J33").PasteSpecial(-4122)
$ws.Range("J33").Value = "node - host"

# --- Row 30's own formula lost its shared-group partner and is now spelled out literally ---
$ws.Range("B30").Formula = "=COUNTA(C28:I28)"

# --- Update the view: scroll back to the top and move the active selection ---
[void]$ws.Range("M33").Select()
